$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NATMI LR-pair table (Lgi4 -> Adam23) re-run with an additional "ECs"
# sending/target cluster, per the diff: rows 2-10 now cover the full 3x3
# {ECs, FAPs, sCs} x {ECs, FAPs, sCs} grid (previously only 2x3, missing ECs
# as a sending cluster).

# Row 2: ECs -> ECs (ligand Lgi4 / receptor Adam23)
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Lgi4"
$ws.Cells.Item(2,3).Value = "Adam23"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.1241926666666667
$ws.Cells.Item(2,8).Value = 0.372578
$ws.Cells.Item(2,9).Value = 0.00757902233016378
$ws.Cells.Item(2,10).Value = 0.00757902233016378
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.1145763333333333
$ws.Cells.Item(2,14).Value = 0.343729
$ws.Cells.Item(2,15).Value = 0.006557053879060051
$ws.Cells.Item(2,16).Value = 0.006557053879060051
$ws.Cells.Item(2,17).Value = 0.01422954037355555
$ws.Cells.Item(2,18).Value = 0.128065863362
$ws.Cells.Item(2,19).Value = 0.00004969605776948316
$ws.Cells.Item(2,20).Value = 0.00004969605776948316

# Row 3: ECs -> FAPs (ligand Lgi4 / receptor Adam23)
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Lgi4"
$ws.Cells.Item(3,3).Value = "Adam23"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.1241926666666667
$ws.Cells.Item(3,8).Value = 0.372578
$ws.Cells.Item(3,9).Value = 0.00757902233016378
$ws.Cells.Item(3,10).Value = 0.00757902233016378
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 9.390663666666667
$ws.Cells.Item(3,14).Value = 28.171991
$ws.Cells.Item(3,15).Value = 0.5374154140831726
$ws.Cells.Item(3,16).Value = 0.5374154140831726
$ws.Cells.Item(3,17).Value = 1.166251562533111
$ws.Cells.Item(3,18).Value = 10.496264062798
$ws.Cells.Item(3,19).Value = 0.00407308342391058
$ws.Cells.Item(3,20).Value = 0.00407308342391058

# Row 4: ECs -> sCs (ligand Lgi4 / receptor Adam23)
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Lgi4"
$ws.Cells.Item(4,3).Value = "Adam23"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.1241926666666667
$ws.Cells.Item(4,8).Value = 0.372578
$ws.Cells.Item(4,9).Value = 0.00757902233016378
$ws.Cells.Item(4,10).Value = 0.00757902233016378
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 7.968512
$ws.Cells.Item(4,14).Value = 23.905536
$ws.Cells.Item(4,15).Value = 0.4560275320377672
$ws.Cells.Item(4,16).Value = 0.4560275320377672
$ws.Cells.Item(4,17).Value = 0.9896307546453332
$ws.Cells.Item(4,18).Value = 8.906676791807998
$ws.Cells.Item(4,19).Value = 0.003456242848483717
$ws.Cells.Item(4,20).Value = 0.003456242848483717

# Row 5: FAPs -> ECs (ligand Lgi4 / receptor Adam23)
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Lgi4"
$ws.Cells.Item(5,3).Value = "Adam23"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 8.566363666666668
$ws.Cells.Item(5,8).Value = 25.699091
$ws.Cells.Item(5,9).Value = 0.5227737133000635
$ws.Cells.Item(5,10).Value = 0.5227737133000636
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.1145763333333333
$ws.Cells.Item(5,14).Value = 0.343729
$ws.Cells.Item(5,15).Value = 0.006557053879060051
$ws.Cells.Item(5,16).Value = 0.006557053879060051
$ws.Cells.Item(5,17).Value = 0.9815025389265556
$ws.Cells.Item(5,18).Value = 8.833522850339
$ws.Cells.Item(5,19).Value = 0.003427855404664808
$ws.Cells.Item(5,20).Value = 0.003427855404664809

# Row 6: FAPs -> FAPs (ligand Lgi4 / receptor Adam23)
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Lgi4"
$ws.Cells.Item(6,3).Value = "Adam23"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 8.566363666666668
$ws.Cells.Item(6,8).Value = 25.699091
$ws.Cells.Item(6,9).Value = 0.5227737133000635
$ws.Cells.Item(6,10).Value = 0.5227737133000636
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 9.390663666666667
$ws.Cells.Item(6,14).Value = 28.171991
$ws.Cells.Item(6,15).Value = 0.5374154140831726
$ws.Cells.Item(6,16).Value = 0.5374154140831726
$ws.Cells.Item(6,17).Value = 80.44384004002012
$ws.Cells.Item(6,18).Value = 723.994560360181
$ws.Cells.Item(6,19).Value = 0.2809466516049514
$ws.Cells.Item(6,20).Value = 0.2809466516049515

# Row 7: FAPs -> sCs (ligand Lgi4 / receptor Adam23)
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Lgi4"
$ws.Cells.Item(7,3).Value = "Adam23"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 8.566363666666668
$ws.Cells.Item(7,8).Value = 25.699091
$ws.Cells.Item(7,9).Value = 0.5227737133000635
$ws.Cells.Item(7,10).Value = 0.5227737133000636
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 7.968512
$ws.Cells.Item(7,14).Value = 23.905536
$ws.Cells.Item(7,15).Value = 0.4560275320377672
$ws.Cells.Item(7,16).Value = 0.4560275320377672
$ws.Cells.Item(7,17).Value = 68.26117167419734
$ws.Cells.Item(7,18).Value = 614.350545067776
$ws.Cells.Item(7,19).Value = 0.2383992062904473
$ws.Cells.Item(7,20).Value = 0.2383992062904473

# Row 8: sCs -> ECs (ligand Lgi4 / receptor Adam23)
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Lgi4"
$ws.Cells.Item(8,3).Value = "Adam23"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 7.695813999999999
$ws.Cells.Item(8,8).Value = 23.087442
$ws.Cells.Item(8,9).Value = 0.4696472643697726
$ws.Cells.Item(8,10).Value = 0.4696472643697726
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.1145763333333333
$ws.Cells.Item(8,14).Value = 0.343729
$ws.Cells.Item(8,15).Value = 0.006557053879060051
$ws.Cells.Item(8,16).Value = 0.006557053879060051
$ws.Cells.Item(8,17).Value = 0.8817581501353332
$ws.Cells.Item(8,18).Value = 7.935823351217999
$ws.Cells.Item(8,19).Value = 0.003079502416625759
$ws.Cells.Item(8,20).Value = 0.003079502416625759

# Row 9: sCs -> FAPs (ligand Lgi4 / receptor Adam23)
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Lgi4"
$ws.Cells.Item(9,3).Value = "Adam23"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 7.695813999999999
$ws.Cells.Item(9,8).Value = 23.087442
$ws.Cells.Item(9,9).Value = 0.4696472643697726
$ws.Cells.Item(9,10).Value = 0.4696472643697726
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 9.390663666666667
$ws.Cells.Item(9,14).Value = 28.171991
$ws.Cells.Item(9,15).Value = 0.5374154140831726
$ws.Cells.Item(9,16).Value = 0.5374154140831726
$ws.Cells.Item(9,17).Value = 72.26880091522466
$ws.Cells.Item(9,18).Value = 650.419208237022
$ws.Cells.Item(9,19).Value = 0.2523956790543106
$ws.Cells.Item(9,20).Value = 0.2523956790543106

# Row 10: sCs -> sCs (ligand Lgi4 / receptor Adam23)
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Lgi4"
$ws.Cells.Item(10,3).Value = "Adam23"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 7.695813999999999
$ws.Cells.Item(10,8).Value = 23.087442
$ws.Cells.Item(10,9).Value = 0.4696472643697726
$ws.Cells.Item(10,10).Value = 0.4696472643697726
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 7.968512
$ws.Cells.Item(10,14).Value = 23.905536
$ws.Cells.Item(10,15).Value = 0.4560275320377672
$ws.Cells.Item(10,16).Value = 0.4560275320377672
$ws.Cells.Item(10,17).Value = 61.32418620876799
$ws.Cells.Item(10,18).Value = 551.917675878912
$ws.Cells.Item(10,19).Value = 0.2141720828988362
$ws.Cells.Item(10,20).Value = 0.2141720828988362

